# Updates the cryptos list with latest price/volume values.
# Generated from commit: "Updated cryptos list on Mon May 27 09:58:59 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value. Values that look like plain numbers are entered
# with a leading apostrophe so Excel keeps them as text (matching the
# original inline-string cell type and preserving formatting such as
# trailing zeros, e.g. "1.00" instead of being coerced to the number 1).
$updates = [ordered]@{
    "D2" = "68.634.06"
    "E2" = "  -0.62%  "
    "D3" = "3.914.30"
    "E3" = "  +3.17%  "
    "E4" = "  +0.09%  "
    "D5" = "'602.78"
    "E5" = "  +0.33%  "
    "D6" = "'165.72"
    "E6" = "  +1.82%  "
    "D7" = "3.912.18"
    "E7" = "  +3.15%  "
    "E8" = "  -0.04%  "
    "D9" = "'0.530"
    "E9" = "  -1.13%  "
    "D10" = "'0.167"
    "E10" = "  -1.40%  "
    "D11" = "'6.39"
    "E11" = "  +1.36%  "
    "D12" = "'0.459"
    "E12" = "  +0.31%  "
    "D13" = "'0.0000259"
    "E13" = "  +5.45%  "
    "D14" = "'37.28"
    "E14" = "  +0.34%  "
    "D15" = "4.570.62"
    "E15" = "  +3.23%  "
    "D16" = "3.901.43"
    "E16" = "  +2.20%  "
    "D17" = "68.749.16"
    "E17" = "  -0.59%  "
    "D18" = "'7.46"
    "E18" = "  +0.96%  "
    "D19" = "'17.11"
    "E19" = "  -1.20%  "
    "E20" = "  -2.15%  "
    "D21" = "'10.98"
    "E21" = "  -2.53%  "
    "D22" = "'486.04"
    "E22" = "  -0.44%  "
    "E23" = "  +0.49%  "
    "D24" = "'0.0000169"
    "E24" = "  +11.19%  "
    "D25" = "'84.48"
    "E25" = "  +0.02%  "
    "E26" = "  -0.49%  "
    "E27" = "  -0.81%  "
    "D28" = "'10.11"
    "E28" = "  +0.71%  "
    "E29" = "  +0.02%  "
    "E30" = "  -0.99%  "
    "D31" = "4.066.73"
    "E31" = "  +3.14%  "
    "D32" = "'2.38"
    "E32" = "  -0.35%  "
    "D33" = "'7.76"
    "E33" = "  -3.43%  "
    "D34" = "'31.85"
    "E34" = "  +0.20%  "
    "D35" = "3.864.81"
    "E35" = "  +3.31%  "
    "E36" = "  +0.35%  "
    "E37" = "  +2.04%  "
    "D38" = "'5.92"
    "E38" = "  +0.64%  "
    "D39" = "'0.137"
    "E39" = "  -1.96%  "
    "D40" = "'3.19"
    "E40" = "  +5.97%  "
    "D41" = "'1.00"
    "E41" = "  +0.05%  "
    "D42" = "'0.315"
    "D43" = "'430.57"
    "E43" = "  +2.98%  "
    "D44" = "'48.44"
    "E44" = "  -0.06%  "
    "E45" = "  +0.14%  "
    "E46" = "  +1.46%  "
    "D48" = "'26.19"
    "E48" = "  +7.63%  "
    "D49" = "'141.81"
    "E49" = "  +0.09%  "
    "D50" = "2.811.38"
    "E50" = "  -0.11%  "
    "D51" = "'0.0353"
    "E51" = "  +0.92%  "
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

